$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "_old" / "_new" header labels (row 1) to the
#    format-version-specific names "_FV2410" / "_FV2504".
$oldToNew = @{
    "Segmentname_old" = "Segmentname_FV2410";
    "Segmentgruppe_old" = "Segmentgruppe_FV2410";
    "Segment_old" = "Segment_FV2410";
    "Datenelement_old" = "Datenelement_FV2410";
    "Segment ID_old" = "Segment ID_FV2410";
    "Code_old" = "Code_FV2410";
    "Qualifier_old" = "Qualifier_FV2410";
    "Beschreibung_old" = "Beschreibung_FV2410";
    "Bedingungsausdruck_old" = "Bedingungsausdruck_FV2410";
    "Bedingung_old" = "Bedingung_FV2410";
    "Segmentname_new" = "Segmentname_FV2504";
    "Segmentgruppe_new" = "Segmentgruppe_FV2504";
    "Segment_new" = "Segment_FV2504";
    "Datenelement_new" = "Datenelement_FV2504";
    "Segment ID_new" = "Segment ID_FV2504";
    "Code_new" = "Code_FV2504";
    "Qualifier_new" = "Qualifier_FV2504";
    "Beschreibung_new" = "Beschreibung_FV2504";
    "Bedingungsausdruck_new" = "Bedingungsausdruck_FV2504";
    "Bedingung_new" = "Bedingung_FV2504";
}

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value()
    if ($oldToNew.ContainsKey($current)) {
        $cell.Value = $oldToNew[$current]
    }
}

# 2. Turn the used range into an Excel table (ListObject) so the header
#    row gets filter buttons and a defined name ("Table1").
$tableRange = $ws.Range("A1:U60")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = $null

# 3. Freeze the header row (split after row 1, top-left of the
#    scrollable area is A2).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
